# Auto-generated: refresh Leve profit-calculation columns (H:N) with new
# market-price snapshot values, per the "update Sheets via scheduled runner" commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 629811.4
$ws.Range("I64").Value = 912607.4399999999
$ws.Range("J64").Value = 7660
$ws.Range("K64").Value = 912607.4399999999
$ws.Range("L64").Value = 7660
$ws.Range("M64").Value = -912359.4399999999
$ws.Range("N64").Value = -8156
$ws.Range("H67").Value = 629811.4
$ws.Range("I67").Value = 912607.4399999999
$ws.Range("J67").Value = 7660
$ws.Range("K67").Value = 912607.4399999999
$ws.Range("L67").Value = 7660
$ws.Range("M67").Value = -911749.4399999999
$ws.Range("N67").Value = -9376
$ws.Range("H74").Value = 3600
$ws.Range("I74").Value = 3457.1428
$ws.Range("J74").Value = 3800
$ws.Range("K74").Value = 3457.1428
$ws.Range("L74").Value = 3800
$ws.Range("M74").Value = -2521.1428
$ws.Range("N74").Value = -5672
$ws.Range("H76").Value = 4833887.5
$ws.Range("I76").Value = 5293886.5
$ws.Range("K76").Value = 5293886.5
$ws.Range("M76").Value = -5293571.5
$ws.Range("H77").Value = 3600
$ws.Range("I77").Value = 3457.1428
$ws.Range("J77").Value = 3800
$ws.Range("K77").Value = 17285.714
$ws.Range("L77").Value = 19000
$ws.Range("M77").Value = -12605.714
$ws.Range("N77").Value = -28360
$ws.Range("H79").Value = 4833887.5
$ws.Range("I79").Value = 5293886.5
$ws.Range("K79").Value = 5293886.5
$ws.Range("M79").Value = -5292794.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18192.111
$ws.Range("I32").Value = 2477.574
$ws.Range("J32").Value = 112479.336
$ws.Range("K32").Value = 2477.574
$ws.Range("L32").Value = 112479.336
$ws.Range("M32").Value = -2190.574
$ws.Range("N32").Value = -113053.336
$ws.Range("H61").Value = 2861.8262
$ws.Range("I61").Value = 2016.2667
$ws.Range("J61").Value = 4447.25
$ws.Range("K61").Value = 2016.2667
$ws.Range("L61").Value = 4447.25
$ws.Range("M61").Value = -1804.2667
$ws.Range("N61").Value = -4871.25
$ws.Range("H136").Value = 2861.8262
$ws.Range("I136").Value = 2016.2667
$ws.Range("J136").Value = 4447.25
$ws.Range("K136").Value = 6048.800099999999
$ws.Range("L136").Value = 13341.75
$ws.Range("M136").Value = -3498.800099999999
$ws.Range("N136").Value = -18441.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1056.45
$ws.Range("I94").Value = 1222.4
$ws.Range("J94").Value = 558.6
$ws.Range("K94").Value = 1222.4
$ws.Range("L94").Value = 558.6
$ws.Range("M94").Value = -771.4000000000001
$ws.Range("N94").Value = -1460.6
$ws.Range("H99").Value = 2143.75
$ws.Range("I99").Value = 2157.1428
$ws.Range("K99").Value = 2157.1428
$ws.Range("M99").Value = -659.1428000000001
$ws.Range("H105").Value = 2748.5
$ws.Range("I105").Value = 2456.1292
$ws.Range("K105").Value = 2456.1292
$ws.Range("M105").Value = -709.1291999999999
$ws.Range("H107").Value = 849.1539
$ws.Range("I107").Value = 894.1
$ws.Range("J107").Value = 699.3333
$ws.Range("K107").Value = 894.1
$ws.Range("L107").Value = 699.3333
$ws.Range("M107").Value = 1025.9
$ws.Range("N107").Value = -4539.3333
$ws.Range("H122").Value = 29985
$ws.Range("J122").Value = 29985
$ws.Range("L122").Value = 29985
$ws.Range("N122").Value = -39785
$ws.Range("H134").Value = 4570.1055
$ws.Range("I134").Value = 2502.8572
$ws.Range("J134").Value = 5776
$ws.Range("K134").Value = 7508.571599999999
$ws.Range("L134").Value = 17328
$ws.Range("M134").Value = -4973.571599999999
$ws.Range("N134").Value = -22398

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1340
$ws.Range("I16").Value = 466.66666
$ws.Range("J16").Value = 2650
$ws.Range("K16").Value = 466.66666
$ws.Range("L16").Value = 2650
$ws.Range("M16").Value = -179.66666
$ws.Range("N16").Value = -3224
$ws.Range("H58").Value = 2643.5334
$ws.Range("I58").Value = 1382.7142
$ws.Range("J58").Value = 3746.75
$ws.Range("K58").Value = 1382.7142
$ws.Range("L58").Value = 3746.75
$ws.Range("M58").Value = -1179.7142
$ws.Range("N58").Value = -4152.75
$ws.Range("H113").Value = 1340
$ws.Range("I113").Value = 466.66666
$ws.Range("J113").Value = 2650
$ws.Range("K113").Value = 466.66666
$ws.Range("L113").Value = 2650
$ws.Range("M113").Value = 1703.33334
$ws.Range("N113").Value = -6990
$ws.Range("H132").Value = 3414.682
$ws.Range("I132").Value = 2880.875
$ws.Range("J132").Value = 4838.1665
$ws.Range("K132").Value = 8642.625
$ws.Range("L132").Value = 14514.4995
$ws.Range("M132").Value = -6112.625
$ws.Range("N132").Value = -19574.4995
$ws.Range("H136").Value = 2643.5334
$ws.Range("I136").Value = 1382.7142
$ws.Range("J136").Value = 3746.75
$ws.Range("K136").Value = 4148.142599999999
$ws.Range("L136").Value = 11240.25
$ws.Range("M136").Value = -1598.142599999999
$ws.Range("N136").Value = -16340.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1887.5
$ws.Range("I109").Value = 875
$ws.Range("J109").Value = 2900
$ws.Range("K109").Value = 2625
$ws.Range("L109").Value = 8700
$ws.Range("M109").Value = -1585
$ws.Range("N109").Value = -10780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2620.8333
$ws.Range("I80").Value = 2600
$ws.Range("J80").Value = 2650
$ws.Range("K80").Value = 2600
$ws.Range("L80").Value = 2650
$ws.Range("M80").Value = -1602
$ws.Range("N80").Value = -4646
$ws.Range("H83").Value = 2620.8333
$ws.Range("I83").Value = 2600
$ws.Range("J83").Value = 2650
$ws.Range("K83").Value = 13000
$ws.Range("L83").Value = 13250
$ws.Range("M83").Value = -8008
$ws.Range("N83").Value = -23234
$ws.Range("H122").Value = 2778980.5
$ws.Range("I122").Value = 3704474
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 11113422
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -11110972
$ws.Range("N122").Value = -12400
$ws.Range("H123").Value = 9070.32
$ws.Range("J123").Value = 9070.32
$ws.Range("L123").Value = 9070.32
$ws.Range("N123").Value = -13970.32

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1154.6666
$ws.Range("I46").Value = 1180
$ws.Range("J46").Value = 1123
$ws.Range("K46").Value = 1180
$ws.Range("L46").Value = 1123
$ws.Range("M46").Value = -992
$ws.Range("N46").Value = -1499
$ws.Range("H61").Value = 4473.5557
$ws.Range("I61").Value = 4157.5415
$ws.Range("J61").Value = 7001.6665
$ws.Range("K61").Value = 4157.5415
$ws.Range("L61").Value = 7001.6665
$ws.Range("M61").Value = -3955.5415
$ws.Range("N61").Value = -7405.6665
$ws.Range("H113").Value = 4473.5557
$ws.Range("I113").Value = 4157.5415
$ws.Range("J113").Value = 7001.6665
$ws.Range("K113").Value = 4157.5415
$ws.Range("L113").Value = 7001.6665
$ws.Range("M113").Value = -1987.5415
$ws.Range("N113").Value = -11341.6665
$ws.Range("H122").Value = 3544.3
$ws.Range("I122").Value = 1721.5
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 5164.5
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2714.5
$ws.Range("N122").Value = -16900
$ws.Range("H136").Value = 4187.1763
$ws.Range("I136").Value = 1993.1538
$ws.Range("J136").Value = 11317.75
$ws.Range("K136").Value = 5979.4614
$ws.Range("L136").Value = 33953.25
$ws.Range("M136").Value = -3429.4614
$ws.Range("N136").Value = -39053.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2778777.5
$ws.Range("I107").Value = 4630606.5
$ws.Range("J107").Value = 1034.375
$ws.Range("K107").Value = 13891819.5
$ws.Range("L107").Value = 3103.125
$ws.Range("M107").Value = -13889899.5
$ws.Range("N107").Value = -6943.125
